# refactored jobs_site class methods to jobs_queue
$wb = $excel.ActiveWorkbook

# --- Sheets ---
$wsJobs = $wb.Worksheets.Item("job_sites")
$wsCity = $wb.Worksheets.Item("city")

# --- Add new job_sites row (id 51 / row 52) ---
# (added before the city row so shared-string indices line up: 199,200,201 then 202)
# Copy formatting down from the last existing row first, then overwrite values.
$wsJobs.Range("A51:J51").Copy()
$wsJobs.Range("A52:J52").PasteSpecial()
$excel.CutCopyMode = $false
$wsJobs.Cells.Item(52, 1).Value = 51
$wsJobs.Cells.Item(52, 2).Value = "hta consulting"
$wsJobs.Cells.Item(52, 3).Value = "http://www.htaconsulting.com/blog/"
$wsJobs.Cells.Item(52, 4).Value = "Research & evaluation firm in bay area. Looks like good work."
$wsJobs.Cells.Item(52, 5).Value = 20
$wsJobs.Cells.Item(52, 6).Value = 0
$wsJobs.Cells.Item(52, 7).Value = 1
$wsJobs.Cells.Item(52, 8).Value = 1
$wsJobs.Cells.Item(52, 9).Value = 1
$wsJobs.Cells.Item(52, 10).Value = 13

# --- Add new city lookup row (Berkeley, index 13) ---
$wsCity.Cells.Item(15, 1).Value = 13
$wsCity.Cells.Item(15, 2).Value = "Berkeley"

# --- Leftover selection on city sheet from editing session ---
$wsCity.Activate()
$wsCity.Range("B39").Select()

# --- View changes on job_sites: freeze top row, move selection to A2 ---
$wsJobs.Activate()
$excel.ActiveWindow.FreezePanes = $false
$wsJobs.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Window screen position tweak (best effort / cosmetic) ---
$excel.ActiveWindow.Left = -28920
$excel.ActiveWindow.Top = 1185
